$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 392; this shifts all rows from 392..470 down to 393..471
$ws.Rows.Item(392).Insert()

# Populate the newly inserted row 392 with the new record (Camote -> Paine, Región Metropolitana -> Región de O'Higgins)
$ws.Range("A392").Value = 11
$ws.Range("B392").Value = "Vega Monumental Concepción"
$ws.Range("C392").Value = "Bíobío"
$ws.Range("D392").Value = 45258
$ws.Range("E392").Value = 8
$ws.Range("F392").Value = 100112045
$ws.Range("G392").Value = "Zapallo"
$ws.Range("H392").Value = "Paine"
$ws.Range("I392").Value = "1a (guarda)"
$ws.Range("J392").Value = 550
$ws.Range("K392").Value = 800
$ws.Range("L392").Value = 1000
$ws.Range("M392").Value = 873
$ws.Range("N392").Value = "$/kilo (volumen en unidades)"
$ws.Range("O392").Value = "Región de O'Higgins"
$ws.Range("P392").Value = 873
$ws.Range("Q392").Value = 1
$ws.Range("R392").Value = "Hortaliza"
